# Loan RBI, Variable Instalments
#
# On the "Repayment schedule" sheet, insert a new blank column before
# column N (the existing "Late"/"heading"/"Outstanding" columns and their
# data all shift one column to the right, from N/O/P -> O/P/Q).
#
# Also update the view state: the "Repayment schedule" sheet becomes the
# active/selected tab (it was "Input" before), with a new selection, and
# the "Input" sheet keeps a plain (non-active) selection too.

$wb = $excel.ActiveWorkbook

# --- Repayment schedule: insert a new column at N ---------------------
$wsSchedule = $wb.Worksheets.Item("Repayment schedule")
$wsSchedule.Columns.Item(14).Insert()

# New column inherits the width of the column immediately to its left (M)
$wsSchedule.Columns.Item(14).ColumnWidth = $wsSchedule.Columns.Item(13).ColumnWidth

# --- Input sheet: update selection, no longer the active tab ----------
$wsInput = $wb.Worksheets.Item("Input")
[void]$wsInput.Range("D22").Select()

# --- Repayment schedule becomes the active tab with a new selection ---
$wsSchedule.Activate()
[void]$wsSchedule.Range("R7").Select()
